# Regenerate the "K" column (column G) values for save_data rows 2-32.
# These are the new strikeout/K counts replacing the prior "Strike#" values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 2
    4  = 2
    5  = 0
    6  = 2
    7  = 2
    8  = 1
    9  = 2
    10 = 2
    11 = 2
    12 = 0
    13 = 0
    14 = 0
    15 = 1
    16 = 2
    17 = 2
    18 = 2
    19 = 2
    20 = 3
    21 = 0
    22 = 5
    23 = 2
    24 = 2
    25 = 1
    26 = 0
    27 = 1
    28 = 1
    29 = 3
    30 = 1
    31 = 2
    32 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
